# Auto-generated: update cryptos price/volume values per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.708.98"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "3.026.56"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "511.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.68%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +1.97%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("E10").Value = "  +2.29%  "
$ws.Range("E11").Value = "  +5.36%  "
$ws.Range("D12").Value = "3.547.84"
$ws.Range("E12").Value = "  +2.48%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("E15").Value = "  +3.95%  "
$ws.Range("D16").Value = "56.707.95"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").Value = "3.028.45"
$ws.Range("E17").Value = "  +2.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "333.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.86%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  +3.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.36%  "
$ws.Range("D25").Value = "3.158.66"
$ws.Range("E25").Value = "  +2.69%  "
$ws.Range("E26").Value = "  +3.24%  "
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("E28").Value = "  +8.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.93%  "
$ws.Range("E33").Value = "  +3.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "152.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "27.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +14.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.62%  "
$ws.Range("E38").Value = "  +2.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0662"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.14%  "
$ws.Range("D40").Value = "3.066.25"
$ws.Range("E40").Value = "  +2.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.98%  "
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("E43").Value = "  +4.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.660"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.64%  "
$ws.Range("D45").Value = "2.213.37"
$ws.Range("E45").Value = "  +3.44%  "
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0243"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.04%  "
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("E51").Value = "  +1.40%  "
